$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1672.8572
$ws.Range("J17").Value = 1672.8572
$ws.Range("L17").Value = 5018.571599999999
$ws.Range("N17").Value = -5354.571599999999
$ws.Range("H33").Value = 314.8
$ws.Range("I33").Value = 286.54544
$ws.Range("J33").Value = 392.5
$ws.Range("K33").Value = 286.54544
$ws.Range("L33").Value = 392.5
$ws.Range("M33").Value = -57.54543999999999
$ws.Range("N33").Value = -850.5
$ws.Range("H38").Value = 858.2308
$ws.Range("I38").Value = 144.625
$ws.Range("J38").Value = 2000
$ws.Range("K38").Value = 433.875
$ws.Range("L38").Value = 6000
$ws.Range("M38").Value = -61.875
$ws.Range("N38").Value = -6744
$ws.Range("H39").Value = 112.125
$ws.Range("I39").Value = 63.545456
$ws.Range("J39").Value = 219
$ws.Range("K39").Value = 190.636368
$ws.Range("L39").Value = 657
$ws.Range("M39").Value = 105.363632
$ws.Range("N39").Value = -1249
$ws.Range("H42").Value = 97.25
$ws.Range("J42").Value = 89
$ws.Range("L42").Value = 267
$ws.Range("N42").Value = -727
$ws.Range("H43").Value = 1770.7084
$ws.Range("I43").Value = 1133.6666
$ws.Range("J43").Value = 1861.7142
$ws.Range("K43").Value = 1133.6666
$ws.Range("L43").Value = 1861.7142
$ws.Range("M43").Value = -1064.6666
$ws.Range("N43").Value = -1999.7142
$ws.Range("H116").Value = 2913.25
$ws.Range("I116").Value = 1780
$ws.Range("J116").Value = 4802
$ws.Range("K116").Value = 1780
$ws.Range("L116").Value = 4802
$ws.Range("M116").Value = 1662
$ws.Range("N116").Value = -11686
$ws.Range("H137").Value = 2443473.2
$ws.Range("I137").Value = 9092485
$ws.Range("J137").Value = 5502.067
$ws.Range("K137").Value = 27277455
$ws.Range("L137").Value = 16506.201
$ws.Range("M137").Value = -27274905
$ws.Range("N137").Value = -21606.201

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1266.95
$ws.Range("I2").Value = 1163.8334
$ws.Range("J2").Value = 1576.3
$ws.Range("K2").Value = 1163.8334
$ws.Range("L2").Value = 1576.3
$ws.Range("M2").Value = -1050.8334
$ws.Range("N2").Value = -1802.3
$ws.Range("H32").Value = 6646.34
$ws.Range("I32").Value = 3665.5903
$ws.Range("J32").Value = 21199.412
$ws.Range("K32").Value = 3665.5903
$ws.Range("L32").Value = 21199.412
$ws.Range("M32").Value = -3378.5903
$ws.Range("N32").Value = -21773.412
$ws.Range("H61").Value = 22268566
$ws.Range("I61").Value = 27806756
$ws.Range("J61").Value = 115805.555
$ws.Range("K61").Value = 27806756
$ws.Range("L61").Value = 115805.555
$ws.Range("M61").Value = -27806544
$ws.Range("N61").Value = -116229.555
$ws.Range("H74").Value = 6708006.5
$ws.Range("I74").Value = 9036933
$ws.Range("J74").Value = 79522.84
$ws.Range("K74").Value = 9036933
$ws.Range("L74").Value = 79522.84
$ws.Range("M74").Value = -9036059
$ws.Range("N74").Value = -81270.84
$ws.Range("H77").Value = 6708006.5
$ws.Range("I77").Value = 9036933
$ws.Range("J77").Value = 79522.84
$ws.Range("K77").Value = 45184665
$ws.Range("L77").Value = 397614.2
$ws.Range("M77").Value = -45180297
$ws.Range("N77").Value = -406350.2
$ws.Range("H114").Value = 25398
$ws.Range("J114").Value = 25398
$ws.Range("L114").Value = 25398
$ws.Range("N114").Value = -34076
$ws.Range("H116").Value = 1266.95
$ws.Range("I116").Value = 1163.8334
$ws.Range("J116").Value = 1576.3
$ws.Range("K116").Value = 1163.8334
$ws.Range("L116").Value = 1576.3
$ws.Range("M116").Value = 1130.1666
$ws.Range("N116").Value = -6164.3
$ws.Range("H122").Value = 15875051
$ws.Range("I122").Value = 1581
$ws.Range("J122").Value = 22224440
$ws.Range("K122").Value = 4743
$ws.Range("L122").Value = 66673320
$ws.Range("M122").Value = -2293
$ws.Range("N122").Value = -66678220
$ws.Range("H132").Value = 52781.562
$ws.Range("I132").Value = 31858.545
$ws.Range("K132").Value = 95575.63499999999
$ws.Range("M132").Value = -93045.63499999999
$ws.Range("H136").Value = 22268566
$ws.Range("I136").Value = 27806756
$ws.Range("J136").Value = 115805.555
$ws.Range("K136").Value = 83420268
$ws.Range("L136").Value = 347416.665
$ws.Range("M136").Value = -83417718
$ws.Range("N136").Value = -352516.665

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1266.95
$ws.Range("I3").Value = 1163.8334
$ws.Range("J3").Value = 1576.3
$ws.Range("K3").Value = 1163.8334
$ws.Range("L3").Value = 1576.3
$ws.Range("M3").Value = -1049.8334
$ws.Range("N3").Value = -1804.3
$ws.Range("H20").Value = 913.4
$ws.Range("I20").Value = 628.2143
$ws.Range("J20").Value = 1276.3636
$ws.Range("K20").Value = 628.2143
$ws.Range("L20").Value = 1276.3636
$ws.Range("M20").Value = -381.2143
$ws.Range("N20").Value = -1770.3636
$ws.Range("H86").Value = 15777.177
$ws.Range("I86").Value = 24301.2
$ws.Range("J86").Value = 3600
$ws.Range("K86").Value = 24301.2
$ws.Range("L86").Value = 3600
$ws.Range("M86").Value = -23178.2
$ws.Range("N86").Value = -5846
$ws.Range("H89").Value = 15777.177
$ws.Range("I89").Value = 24301.2
$ws.Range("J89").Value = 3600
$ws.Range("K89").Value = 121506
$ws.Range("L89").Value = 18000
$ws.Range("M89").Value = -115890
$ws.Range("N89").Value = -29232
$ws.Range("H107").Value = 3365.7827
$ws.Range("I107").Value = 2888.1428
$ws.Range("J107").Value = 4108.778
$ws.Range("K107").Value = 2888.1428
$ws.Range("L107").Value = 4108.778
$ws.Range("M107").Value = -968.1428000000001
$ws.Range("N107").Value = -7948.778
$ws.Range("H134").Value = 1490.7838
$ws.Range("I134").Value = 1411
$ws.Range("J134").Value = 2254.4285
$ws.Range("K134").Value = 4233
$ws.Range("L134").Value = 6763.2855
$ws.Range("M134").Value = -1698
$ws.Range("N134").Value = -11833.2855

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 36024.5
$ws.Range("I132").Value = 2372.524
$ws.Range("J132").Value = 114545.78
$ws.Range("K132").Value = 7117.572
$ws.Range("L132").Value = 343637.34
$ws.Range("M132").Value = -4587.572
$ws.Range("N132").Value = -348697.34

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1179.3158
$ws.Range("I122").Value = 420
$ws.Range("K122").Value = 3780
$ws.Range("M122").Value = -1330
$ws.Range("H131").Value = 985.1799999999999
$ws.Range("J131").Value = 1001.2292
$ws.Range("L131").Value = 3003.6876
$ws.Range("N131").Value = -13083.6876
$ws.Range("H136").Value = 3421.0715
$ws.Range("I136").Value = 3105.8
$ws.Range("J136").Value = 3596.2222
$ws.Range("K136").Value = 9317.400000000001
$ws.Range("L136").Value = 10788.6666
$ws.Range("M136").Value = -4217.400000000001
$ws.Range("N136").Value = -20988.6666

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1808.7368
$ws.Range("I122").Value = 1826.9231
$ws.Range("K122").Value = 5480.7693
$ws.Range("M122").Value = -3030.7693
$ws.Range("H126").Value = 1360.8889
$ws.Range("I126").Value = 978.2857
$ws.Range("K126").Value = 2934.8571
$ws.Range("M126").Value = -464.8571000000002
$ws.Range("H132").Value = 101710.75
$ws.Range("I132").Value = 63732.625
$ws.Range("J132").Value = 253623.25
$ws.Range("K132").Value = 191197.875
$ws.Range("L132").Value = 760869.75
$ws.Range("M132").Value = -188667.875
$ws.Range("N132").Value = -765929.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6912.645
$ws.Range("I22").Value = 8386
$ws.Range("J22").Value = 773.6667
$ws.Range("K22").Value = 8386
$ws.Range("L22").Value = 773.6667
$ws.Range("M22").Value = -8091
$ws.Range("N22").Value = -1363.6667
$ws.Range("H27").Value = 6912.645
$ws.Range("I27").Value = 8386
$ws.Range("J27").Value = 773.6667
$ws.Range("K27").Value = 8386
$ws.Range("L27").Value = 773.6667
$ws.Range("M27").Value = -8279
$ws.Range("N27").Value = -987.6667
$ws.Range("H46").Value = 1783157.8
$ws.Range("I46").Value = 2331622.5
$ws.Range("J46").Value = 647.5
$ws.Range("K46").Value = 2331622.5
$ws.Range("L46").Value = 647.5
$ws.Range("M46").Value = -2331434.5
$ws.Range("N46").Value = -1023.5
$ws.Range("H61").Value = 1827.6428
$ws.Range("I61").Value = 1473.9166
$ws.Range("J61").Value = 3950
$ws.Range("K61").Value = 1473.9166
$ws.Range("L61").Value = 3950
$ws.Range("M61").Value = -1271.9166
$ws.Range("N61").Value = -4354
$ws.Range("H82").Value = 1397.5834
$ws.Range("I82").Value = 1053
$ws.Range("J82").Value = 1880
$ws.Range("K82").Value = 1053
$ws.Range("L82").Value = 1880
$ws.Range("M82").Value = -692
$ws.Range("N82").Value = -2602
$ws.Range("H85").Value = 1397.5834
$ws.Range("I85").Value = 1053
$ws.Range("J85").Value = 1880
$ws.Range("K85").Value = 1053
$ws.Range("L85").Value = 1880
$ws.Range("M85").Value = 195
$ws.Range("N85").Value = -4376
$ws.Range("H93").Value = 1754.8182
$ws.Range("I93").Value = 1577.6666
$ws.Range("K93").Value = 1577.6666
$ws.Range("M93").Value = -329.6666
$ws.Range("H100").Value = 1361
$ws.Range("I100").Value = 1094.9333
$ws.Range("J100").Value = 1804.4445
$ws.Range("K100").Value = 1094.9333
$ws.Range("L100").Value = 1804.4445
$ws.Range("M100").Value = -553.9332999999999
$ws.Range("N100").Value = -2886.4445
$ws.Range("H113").Value = 1827.6428
$ws.Range("I113").Value = 1473.9166
$ws.Range("J113").Value = 3950
$ws.Range("K113").Value = 1473.9166
$ws.Range("L113").Value = 3950
$ws.Range("M113").Value = 696.0834
$ws.Range("N113").Value = -8290
$ws.Range("H132").Value = 102595.3
$ws.Range("I132").Value = 1989
$ws.Range("J132").Value = 169666.17
$ws.Range("K132").Value = 5967
$ws.Range("L132").Value = 508998.51
$ws.Range("M132").Value = -3437
$ws.Range("N132").Value = -514058.51
$ws.Range("H136").Value = 59616.805
$ws.Range("I136").Value = 38313.965
$ws.Range("J136").Value = 147871.42
$ws.Range("K136").Value = 114941.895
$ws.Range("L136").Value = 443614.26
$ws.Range("M136").Value = -112391.895
$ws.Range("N136").Value = -448714.26

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 43516.812
$ws.Range("I132").Value = 26525.3
$ws.Range("J132").Value = 128474.375
$ws.Range("K132").Value = 79575.89999999999
$ws.Range("L132").Value = 385423.125
$ws.Range("M132").Value = -77045.89999999999
$ws.Range("N132").Value = -390483.125
$ws.Range("H136").Value = 52350.977
$ws.Range("I136").Value = 36829.785
$ws.Range("J136").Value = 88567.086
$ws.Range("K136").Value = 110489.355
$ws.Range("L136").Value = 265701.258
$ws.Range("M136").Value = -107939.355
$ws.Range("N136").Value = -270801.258

Write-Host "Applied all cell updates"